# Fruta / hortaliza, semanal
# The weekly consolidation re-shuffled the data rows (2-34) of the sheet.
# Every destination row ends up with exactly the same full set of column
# values (A:T) that some specific source row used to have - i.e. the data
# rows were permuted. We snapshot the old rows first (so we don't
# overwrite data we still need to read), then write the rows back out in
# their new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> old row number (where its data used to live)
$mapping = @{
    2  = 34
    3  = 11
    4  = 12
    5  = 15
    6  = 16
    7  = 9
    8  = 10
    9  = 32
    10 = 27
    11 = 4
    12 = 5
    13 = 6
    14 = 28
    15 = 13
    16 = 14
    17 = 22
    18 = 23
    19 = 24
    20 = 7
    21 = 8
    22 = 31
    23 = 26
    24 = 17
    25 = 18
    26 = 19
    27 = 20
    28 = 3
    29 = 25
    30 = 33
    31 = 2
    32 = 29
    33 = 30
    34 = 21
}

# Snapshot every data row (columns A:T) before any writes happen.
$snapshot = @{}
for ($r = 2; $r -le 34; $r++) {
    $rowRange = $ws.Range("A" + $r + ":T" + $r)
    $snapshot[$r] = $rowRange.Value()
}

# Now write each destination row using the snapshot of its source row.
for ($r = 2; $r -le 34; $r++) {
    $srcRow = $mapping[$r]
    $destRange = $ws.Range("A" + $r + ":T" + $r)
    $destRange.Value = $snapshot[$srcRow]
}
